$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number that was bumped by one
# day (2023-09-10 -> 2023-09-11, i.e. 45179 -> 45180) for every data row.
# Data rows run from row 2 through row 119.
$ws.Range("C2:C119").Value = 45180
